$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-16 11:06:44"

$wsZhCn.Range("H2").Value = "2016-08-16 11:06:37"
$wsZhCn.Range("K2").Value = "2016-08-16 11:07:12"

$wsDeDe.Range("H2").Value = "2016-08-16 11:06:44"
$wsDeDe.Range("K2").Value = "2016-08-16 11:07:20"
